# Add NSW additional Anzac Day observance (2026+)
#
# When Anzac Day (25 Apr) falls on a Saturday or a Sunday, NSW grants an
# additional public holiday on the following Monday:
#   - 2026: Anzac Day is Saturday 25 Apr -> additional Monday 27 Apr
#   - 2027: Anzac Day is Sunday 25 Apr   -> additional Monday 26 Apr
#
# Insert a new row directly below each "Anzac Day" row for those years
# containing the "Anzac Day (additional)" entry, pushing every later row
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-HolidayRow {
    param(
        [int]$RowIndex,
        [string]$DateText,
        [string]$Name
    )

    # Insert a blank row, shifting the current row (and everything below)
    # down by one.
    $ws.Rows.Item($RowIndex).Insert()

    # Force the date column to be stored as text (matching the rest of the
    # sheet) instead of letting Excel auto-convert the "yyyy-mm-dd" string
    # into a date serial number.
    $dateCell = $ws.Range("A" + $RowIndex)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $DateText

    $ws.Range("B" + $RowIndex).Value = $Name
}

# Row 8 holds "2026-04-25 / Anzac Day" -> insert the additional holiday
# right after it, at row 9.
Insert-HolidayRow 9 "2026-04-27" "Anzac Day (additional)"

# After the previous insert, "2027-04-25 / Anzac Day" (originally row 21)
# is now at row 22 -> insert the additional holiday right after it, at row 23.
Insert-HolidayRow 23 "2027-04-26" "Anzac Day (additional)"
